# "caracteristicas dispositivos capa 2"
# Inserts a new worksheet "cisco-model" between "switching" and "AND",
# fills it with a characteristics/description table (Tabla4), adds a small
# "Comandos" lookup column (D6:D10) on the "switching" sheet, and leaves
# the "AND" sheet's data untouched (it just shifts from position 2 to 3).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "switching"
$ws3 = $wb.Worksheets.Item(2)   # "AND" (will end up 3rd after insert)

# --- 1. Update the "switching" sheet's selection/scroll BEFORE inserting
#        the new sheet, while it is still the active one, so the final
#        saved selection on this sheet is D10 (scrolled to show row 5+).
$ws1.Activate()
$ws1.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

# --- 2. Insert the new worksheet "cisco-model" right after "switching"
#        (i.e. before "AND"). Newly added sheets become the active sheet,
#        which matches the target's activeTab/tabSelected state.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "cisco-model"

# --- 3. Populate "cisco-model" column A (rows 2-5) first, matching the
#        shared-string table order of the original authored edit.
$newSheet.Range("A2").Value = "Escalabilidad"
$newSheet.Range("A3").Value = "Tolerancia"
$newSheet.Range("A4").Value = "Disponibilidad"
$newSheet.Range("A5").Value = "Costos"

# --- 4. Add the small "Comandos" reference values in column D of
#        "switching", rows 6-10. NOTE: the insertion order below matches
#        the shared-string table order of the original authored edit
#        (Comandos, DUPLEX, MDIX, SPEED, SWITCHPORT) even though the
#        values land in rows 6,8,7,9,10 respectively.
$ws1.Range("D6").Value  = "Comandos"
$ws1.Range("D8").Value  = "DUPLEX"
$ws1.Range("D7").Value  = "MDIX"
$ws1.Range("D9").Value  = "SPEED"
$ws1.Range("D10").Value = "SWITCHPORT"

# --- 5. Continue populating "cisco-model" in the exact entry order that
#        reproduces the target shared-string sequence.
$newSheet.Range("B2").Value = "Capacidad de Expansion de la topologia"
$newSheet.Range("B3").Value = "Capacidad para reponerse ante un fallo"
$newSheet.Range("A6").Value = "Calidad (QoS)"
$newSheet.Range("A7").Value = "Seguridad"
$newSheet.Range("B4").Value = "Capacidad de poder establecer comunicaciones"
$newSheet.Range("B5").Value = "Capacidad para identificar los mejores trayectos"
$newSheet.Range("B6").Value = "Capacidad para distribuir el ancho de banda"
$newSheet.Range("B7").Value = "Capacidad para prever posibles ataques"
$newSheet.Range("A1").Value = "Caracteristica"
$newSheet.Range("B1").Value = "Descripcion"

# --- 6. Column widths approximating the authored layout.
$newSheet.Columns.Item(1).ColumnWidth = 15.6
$newSheet.Columns.Item(2).ColumnWidth = 43.17

# --- 7. Turn the A1:B7 range into a table ("Tabla4") styled like the
#        workbook's other tables.
$tbl = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:B7"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Tabla4"
$tbl.TableStyle = "TableStyleMedium16"

# --- 8. Make sure "cisco-model" ends up the active/selected sheet & cell,
#        matching the target's tabSelected="1" / activeTab="1".
$newSheet.Activate()
$newSheet.Range("B2").Select()
